# Refresh the cryptos table (price + 1h volume columns) with the latest
# scraped values, as produced by the "Updated cryptos list ... with GitHub
# Actions" automation run. Row 37/38 additionally swap places because the
# source ranking re-ordered RenderToken and TheSandbox.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range('D2').Value = '27.245.55'
$ws.Range('E2').Value = '  +0.62%  '

# Row 3: Ethereum
$ws.Range('D3').Value = '1.905.47'
$ws.Range('E3').Value = '  +0.64%  '

# Row 4: TetherUSD
$ws.Range('E4').Value = '  -0.17%  '

# Row 5: BNB
$ws.Range('D5').Value = '''306.32'
$ws.Range('E5').Value = '  -0.14%  '

# Row 7: XRP
$ws.Range('D7').Value = '''0.5414'
$ws.Range('E7').Value = '  +3.60%  '

# Row 8: Cardano
$ws.Range('E8').Value = '  +1.49%  '

# Row 9: Dogecoin
$ws.Range('D9').Value = '''0.07295'
$ws.Range('E9').Value = '  +0.42%  '

# Row 10: Solana
$ws.Range('D10').Value = '''22.15'
$ws.Range('E10').Value = '  +4.98%  '

# Row 11: Polygon
$ws.Range('D11').Value = '''0.9043'
$ws.Range('E11').Value = '  +0.54%  '

# Row 12: TRON
$ws.Range('D12').Value = '''0.08200'
$ws.Range('E12').Value = '  +0.30%  '

# Row 13: Litecoin
$ws.Range('D13').Value = '''95.95'
$ws.Range('E13').Value = '  -0.38%  '

# Row 14: Polkadot
$ws.Range('D14').Value = '''5.354'
$ws.Range('E14').Value = '  +1.17%  '

# Row 15: BinanceUSD
$ws.Range('D15').Value = '''1.000'
$ws.Range('E15').Value = '  -0.19%  '

# Row 16: Avalanche
$ws.Range('D16').Value = '''14.89'
$ws.Range('E16').Value = '  +2.21%  '

# Row 17: ShibaInu
$ws.Range('D17').Value = '''0.000008650'
$ws.Range('E17').Value = '  +0.91%  '

# Row 18: Dai
$ws.Range('E18').Value = '  -0.32%  '

# Row 19: WrappedEther
$ws.Range('D19').Value = '1.236.87'
$ws.Range('E19').Value = '  -35.99%  '

# Row 20: WrappedBTC
$ws.Range('D20').Value = '27.276.49'
$ws.Range('E20').Value = '  +0.61%  '

# Row 21: Uniswap
$ws.Range('E21').Value = '  -0.45%  '

# Row 22: Cosmos
$ws.Range('D22').Value = '''10.83'

# Row 23: Chainlink
$ws.Range('D23').Value = '''6.521'
$ws.Range('E23').Value = '  +1.78%  '

# Row 24: Monero
$ws.Range('D24').Value = '''148.44'
$ws.Range('E24').Value = '  -0.15%  '

# Row 25: LidoDAOToken
$ws.Range('D25').Value = '''2.304'
$ws.Range('E25').Value = '  +0.64%  '

# Row 26: EthereumClassic
$ws.Range('D26').Value = '''18.39'
$ws.Range('E26').Value = '  +1.17%  '

# Row 27: Toncoin
$ws.Range('D27').Value = '''1.756'
$ws.Range('E27').Value = '  +1.40%  '

# Row 28: BitcoinCash
$ws.Range('D28').Value = '''116.93'
$ws.Range('E28').Value = '  +1.74%  '

# Row 29: InternetComputer(DFINITY)
$ws.Range('D29').Value = '''4.858'
$ws.Range('E29').Value = '  +1.55%  '

# Row 30: Filecoin
$ws.Range('D30').Value = '''4.675'
$ws.Range('E30').Value = '  -3.31%  '

# Row 31: Stellar
$ws.Range('D31').Value = '''0.09207'
$ws.Range('E31').Value = '  -0.18%  '

# Row 32: ImmutableX
$ws.Range('D32').Value = '''0.8277'
$ws.Range('E32').Value = '  +4.95%  '

# Row 33: Hedera
$ws.Range('D33').Value = '''0.05077'
$ws.Range('E33').Value = '  +0.82%  '

# Row 34: ARBITRUM
$ws.Range('D34').Value = '''1.223'
$ws.Range('E34').Value = '  +1.07%  '

# Row 35: HuobiToken
$ws.Range('E35').Value = '  +0.92%  '

# Row 36: MXToken
$ws.Range('D36').Value = '''3.322'
$ws.Range('E36').Value = '  -3.26%  '

# Row 37: RenderToken
$ws.Range('B37').Value = 'TheSandbox'
$ws.Range('C37').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D37').Value = '''0.6072'
$ws.Range('E37').Value = '  +6.11%  '

# Row 38: TheSandbox
$ws.Range('B38').Value = 'RenderToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D38').Value = '''2.696'
$ws.Range('E38').Value = '  +3.66%  '

# Row 39: VeChain
$ws.Range('D39').Value = '''0.02000'
$ws.Range('E39').Value = '  +0.81%  '

# Row 40: TrustWalletToken
$ws.Range('D40').Value = '''1.075'
$ws.Range('E40').Value = '  +0.14%  '

# Row 41: Aptos
$ws.Range('D41').Value = '''9.296'
$ws.Range('E41').Value = '  +2.85%  '

# Row 42: FraxShare
$ws.Range('D42').Value = '''6.668'
$ws.Range('E42').Value = '  +1.79%  '

# Row 43: Quant
$ws.Range('D43').Value = '''116.39'
$ws.Range('E43').Value = '  +0.12%  '

# Row 44: Decentraland
$ws.Range('D44').Value = '''0.5185'
$ws.Range('E44').Value = '  +6.64%  '

# Row 45: Algorand
$ws.Range('D45').Value = '''0.1535'
$ws.Range('E45').Value = '  +1.26%  '

# Row 46: EnergySwap
$ws.Range('D46').Value = '''10.19'
$ws.Range('E46').Value = '  +1.48%  '

# Row 47: PaxDollar
$ws.Range('D47').Value = '''0.9991'
$ws.Range('E47').Value = '  -0.31%  '

# Row 48: NEARProtocol
$ws.Range('D48').Value = '''1.645'
$ws.Range('E48').Value = '  +1.29%  '

# Row 49: Elrond
$ws.Range('D49').Value = '''38.26'
$ws.Range('E49').Value = '  +0.33%  '

# Row 50: Cronos
$ws.Range('D50').Value = '''0.06101'
$ws.Range('E50').Value = '  +2.92%  '

# Row 51: Aave
$ws.Range('D51').Value = '''63.77'
$ws.Range('E51').Value = '  +0.44%  '
